# REC07. tweaks to D64TM & D6434TC
#  - D64TM: corpse removal for test facility after arena reset (row 148 qty change)
#  - D6434TC: revised unmaker secret door (new LINE SPECIALS entries)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D64TM: Corpse removal for test facility after arena reset ---
# Row 148 quantity edited from 2024 to 556
$ws.Range("A148").Value = 556

# --- D6434TC: Revised unmaker secret door ---
# New "TC THINGS" entries (rows 159-161) and a new "LINE SPECIALS" section
# (rows 164-165) appended below the existing table.

# Row 159: Orange Demon artefact
$ws.Range("A159").Value = 4000
$ws.Range("C159").Value = "Orange Demon artefact"

# Row 165: Secret Line Cross
$ws.Range("A165").Value = 994
$ws.Range("A165").HorizontalAlignment = -4131
$ws.Range("C165").Value = "Secret Line Cross"

# Row 164: new section header "LINE SPECIALS" (copy formatting from the
# existing "TC THINGS" header at A144 so it reuses the bold+underline style)
$ws.Range("A144").Copy() | Out-Null
$ws.Range("A164").PasteSpecial(-4122) | Out-Null
$ws.Range("A164").Value = "LINE SPECIALS"
$excel.CutCopyMode = 0

# Row 160: Stimpack
$ws.Range("A160").Value = 7582
$ws.Range("C160").Value = "Stimpack"

# Row 161: Megasphere
$ws.Range("A161").Value = 554
$ws.Range("C161").Value = "Megasphere"

# Update the selection/active cell to match the author's final view state
$ws.Range("C161").Select() | Out-Null
